$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple text updates (date + one corrected room number) ---
$ws.Range("A1").Value = "27/10/2019"
$ws.Range("D13").Value = "AD B310 admin B3"

# --- Fill previously-empty cells with new schedule entries ---
# group cells by the template style they should copy, then set their text
$fillMap_wrap48 = [ordered]@{
    "J13" = "WP-II B314 teacher4 B1"
    "K13" = "OS B317 teacher1 B2"
    "L13" = "OOSE B315 teacher2 B3"
    "M13" = "AD B313 admin B4"
    "R13" = "CG B310 teacher3 B1"
    "S13" = "WP-II B313 teacher4 B2"
    "U13" = "OOSE B315 teacher2 B4"
    "N16" = "WP-II B314 teacher4 B1"
    "O16" = "WP-II B313 teacher4 B2"
    "P16" = "WP-II B312 teacher4 B3"
}
foreach ($cell in $fillMap_wrap48.Keys) {
    $ws.Range("B13").Copy()
    $ws.Range($cell).PasteSpecial(-4122)
    $ws.Range($cell).Value = $fillMap_wrap48[$cell]
}
$fillMap_nowrap_center = [ordered]@{
    "N13" = "CG B305 teacher3"
    "N14" = "OOSE B305 teacher2"
    "J16" = "OS B305 teacher1"
    "R16" = "CG B305 teacher3"
    "J17" = "AD B305 admin"
    "R17" = "OOSE B305 teacher2"
    "J18" = "OOSE B305 teacher2"
    "N18" = "ITC B305 teacher2"
    "R18" = "ITC B305 teacher2"
    "N19" = "AD B305 admin"
}
foreach ($cell in $fillMap_nowrap_center.Keys) {
    $ws.Range("B16").Copy()
    $ws.Range($cell).PasteSpecial(-4122)
    $ws.Range($cell).Value = $fillMap_nowrap_center[$cell]
}
$fillMap_wrap52 = [ordered]@{
    "T13" = "OS B317 teacher1 B3"
    "Q16" = "WP-II B315 teacher4 B4"
}
foreach ($cell in $fillMap_wrap52.Keys) {
    $ws.Range("G16").Copy()
    $ws.Range($cell).PasteSpecial(-4122)
    $ws.Range($cell).Value = $fillMap_wrap52[$cell]
}

# --- Re-style cells that remain empty but change visual style (box -> blank) ---
$blankCells = @("O13", "P13", "Q13", "B14", "C14", "D14", "E14", "F14", "G14", "H14", "I14", "J14", "K14", "L14", "M14", "O14", "P14", "Q14", "R14", "S14", "T14", "U14", "K16", "L16", "M16", "S16", "T16", "U16", "K17", "L17", "M17", "N17", "O17", "P17", "Q17", "S17", "T17", "U17", "K18", "L18", "M18", "O18", "P18", "Q18", "S18", "T18", "U18", "O19", "P19", "Q19")
foreach ($cell in $blankCells) {
    $ws.Range("C16").Copy()
    $ws.Range($cell).PasteSpecial(-4122)
}

# --- Merge cells to match new layout ---
$newMerges = @("J13:J14", "K13:K14", "L13:L14", "M13:M14", "J16:M16", "J17:M17", "J18:M18", "N13:Q13", "N14:Q14", "N16:N17", "O16:O17", "P16:P17", "Q16:Q17", "N18:Q18", "N19:Q19", "R13:R14", "S13:S14", "T13:T14", "U13:U14", "R16:U16", "R17:U17", "R18:U18")
foreach ($m in $newMerges) {
    $ws.Range($m).Merge()
}

Write-Host "Edit complete"